$wb = $excel.ActiveWorkbook

# 1. Rename the "Requested quantity" header on the existing sheets
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2. Add the new "PO Forecast" sheet after "Monthly Trend"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# 3. Header row
$wsForecast.Cells.Item(1,1).Value = "ds"
$wsForecast.Cells.Item(1,2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value = "yhat_upper"

# 4. Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$wsForecast.Cells.Item(2,1).Value = 45011.99999999999
$wsForecast.Cells.Item(2,2).Value = 29
$wsForecast.Cells.Item(2,3).Value = -26.20248491409725
$wsForecast.Cells.Item(2,4).Value = 85.4761265210271
$wsForecast.Cells.Item(3,1).Value = 45025.99999999999
$wsForecast.Cells.Item(3,2).Value = 30
$wsForecast.Cells.Item(3,3).Value = -32.40117467225442
$wsForecast.Cells.Item(3,4).Value = 81.62268289246235
$wsForecast.Cells.Item(4,1).Value = 45032.99999999999
$wsForecast.Cells.Item(4,2).Value = 30
$wsForecast.Cells.Item(4,3).Value = -29.00874033777513
$wsForecast.Cells.Item(4,4).Value = 87.67397200840361
$wsForecast.Cells.Item(5,1).Value = 45046.99999999999
$wsForecast.Cells.Item(5,2).Value = 30
$wsForecast.Cells.Item(5,3).Value = -26.71572849790796
$wsForecast.Cells.Item(5,4).Value = 82.37818534629997
$wsForecast.Cells.Item(6,1).Value = 45053.99999999999
$wsForecast.Cells.Item(6,2).Value = 31
$wsForecast.Cells.Item(6,3).Value = -27.25062417580597
$wsForecast.Cells.Item(6,4).Value = 89.0227325045131
$wsForecast.Cells.Item(7,1).Value = 45060.99999999999
$wsForecast.Cells.Item(7,2).Value = 31
$wsForecast.Cells.Item(7,3).Value = -25.29148418170701
$wsForecast.Cells.Item(7,4).Value = 87.24641683208256
$wsForecast.Cells.Item(8,1).Value = 45067.99999999999
$wsForecast.Cells.Item(8,2).Value = 31
$wsForecast.Cells.Item(8,3).Value = -27.20742502860049
$wsForecast.Cells.Item(8,4).Value = 87.01764344492274
$wsForecast.Cells.Item(9,1).Value = 45074.99999999999
$wsForecast.Cells.Item(9,2).Value = 31
$wsForecast.Cells.Item(9,3).Value = -22.81788339447456
$wsForecast.Cells.Item(9,4).Value = 91.27839488468646
$wsForecast.Cells.Item(10,1).Value = 45081.99999999999
$wsForecast.Cells.Item(10,2).Value = 32
$wsForecast.Cells.Item(10,3).Value = -27.10316124084338
$wsForecast.Cells.Item(10,4).Value = 90.88921342433615
$wsForecast.Cells.Item(11,1).Value = 45088.99999999999
$wsForecast.Cells.Item(11,2).Value = 32
$wsForecast.Cells.Item(11,3).Value = -24.76124807233443
$wsForecast.Cells.Item(11,4).Value = 91.64615605063599
$wsForecast.Cells.Item(12,1).Value = 45095.99999999999
$wsForecast.Cells.Item(12,2).Value = 32
$wsForecast.Cells.Item(12,3).Value = -21.57688511695857
$wsForecast.Cells.Item(12,4).Value = 87.59971236614753
$wsForecast.Cells.Item(13,1).Value = 45102.99999999999
$wsForecast.Cells.Item(13,2).Value = 32
$wsForecast.Cells.Item(13,3).Value = -21.4504721158102
$wsForecast.Cells.Item(13,4).Value = 86.57631429161236
$wsForecast.Cells.Item(14,1).Value = 45109.99999999999
$wsForecast.Cells.Item(14,2).Value = 33
$wsForecast.Cells.Item(14,3).Value = -22.31835780180636
$wsForecast.Cells.Item(14,4).Value = 88.85787721797227
$wsForecast.Cells.Item(15,1).Value = 45116.99999999999
$wsForecast.Cells.Item(15,2).Value = 33
$wsForecast.Cells.Item(15,3).Value = -23.55442308785434
$wsForecast.Cells.Item(15,4).Value = 86.97218636144002
$wsForecast.Cells.Item(16,1).Value = 45123.99999999999
$wsForecast.Cells.Item(16,2).Value = 33
$wsForecast.Cells.Item(16,3).Value = -26.11680070416986
$wsForecast.Cells.Item(16,4).Value = 91.78442767617126
$wsForecast.Cells.Item(17,1).Value = 45130.99999999999
$wsForecast.Cells.Item(17,2).Value = 33
$wsForecast.Cells.Item(17,3).Value = -30.23040350855943
$wsForecast.Cells.Item(17,4).Value = 86.28770533455818
$wsForecast.Cells.Item(18,1).Value = 45137.99999999999
$wsForecast.Cells.Item(18,2).Value = 34
$wsForecast.Cells.Item(18,3).Value = -24.27400506336808
$wsForecast.Cells.Item(18,4).Value = 89.26855638314358
$wsForecast.Cells.Item(19,1).Value = 45144.99999999999
$wsForecast.Cells.Item(19,2).Value = 34
$wsForecast.Cells.Item(19,3).Value = -25.68596670910047
$wsForecast.Cells.Item(19,4).Value = 89.19667971396554
$wsForecast.Cells.Item(20,1).Value = 45151.99999999999
$wsForecast.Cells.Item(20,2).Value = 34
$wsForecast.Cells.Item(20,3).Value = -23.41455695860182
$wsForecast.Cells.Item(20,4).Value = 90.23858289676231
$wsForecast.Cells.Item(21,1).Value = 45158.99999999999
$wsForecast.Cells.Item(21,2).Value = 34
$wsForecast.Cells.Item(21,3).Value = -23.79035375576454
$wsForecast.Cells.Item(21,4).Value = 95.25881651663148
$wsForecast.Cells.Item(22,1).Value = 45165.99999999999
$wsForecast.Cells.Item(22,2).Value = 35
$wsForecast.Cells.Item(22,3).Value = -20.38971938690293
$wsForecast.Cells.Item(22,4).Value = 88.46503106736412
$wsForecast.Cells.Item(23,1).Value = 45172.99999999999
$wsForecast.Cells.Item(23,2).Value = 35
$wsForecast.Cells.Item(23,3).Value = -21.21678413411461
$wsForecast.Cells.Item(23,4).Value = 90.07487139235923
$wsForecast.Cells.Item(24,1).Value = 45186.99999999999
$wsForecast.Cells.Item(24,2).Value = 35
$wsForecast.Cells.Item(24,3).Value = -23.433250732927
$wsForecast.Cells.Item(24,4).Value = 90.38831545893964
$wsForecast.Cells.Item(25,1).Value = 45193.99999999999
$wsForecast.Cells.Item(25,2).Value = 36
$wsForecast.Cells.Item(25,3).Value = -23.15149524930697
$wsForecast.Cells.Item(25,4).Value = 88.2262577905729
$wsForecast.Cells.Item(26,1).Value = 45200.99999999999
$wsForecast.Cells.Item(26,2).Value = 36
$wsForecast.Cells.Item(26,3).Value = -22.17375633749975
$wsForecast.Cells.Item(26,4).Value = 93.2057215820582
$wsForecast.Cells.Item(27,1).Value = 45221.99999999999
$wsForecast.Cells.Item(27,2).Value = 37
$wsForecast.Cells.Item(27,3).Value = -17.70505007320202
$wsForecast.Cells.Item(27,4).Value = 95.61693375092878
$wsForecast.Cells.Item(28,1).Value = 45277.99999999999
$wsForecast.Cells.Item(28,2).Value = 39
$wsForecast.Cells.Item(28,3).Value = -16.49852438236568
$wsForecast.Cells.Item(28,4).Value = 96.289317059458
$wsForecast.Cells.Item(29,1).Value = 45284.99999999999
$wsForecast.Cells.Item(29,2).Value = 39
$wsForecast.Cells.Item(29,3).Value = -15.34942060491667
$wsForecast.Cells.Item(29,4).Value = 98.37311416780535
$wsForecast.Cells.Item(30,1).Value = 45291.99999999999
$wsForecast.Cells.Item(30,2).Value = 39
$wsForecast.Cells.Item(30,3).Value = -19.88289009849426
$wsForecast.Cells.Item(30,4).Value = 89.45941897710379
$wsForecast.Cells.Item(31,1).Value = 45298.99999999999
$wsForecast.Cells.Item(31,2).Value = 39
$wsForecast.Cells.Item(31,3).Value = -16.85218641083066
$wsForecast.Cells.Item(31,4).Value = 96.66309438425577
$wsForecast.Cells.Item(32,1).Value = 45305.99999999999
$wsForecast.Cells.Item(32,2).Value = 40
$wsForecast.Cells.Item(32,3).Value = -18.87798540534655
$wsForecast.Cells.Item(32,4).Value = 97.51723219421991
$wsForecast.Cells.Item(33,1).Value = 45312.99999999999
$wsForecast.Cells.Item(33,2).Value = 40
$wsForecast.Cells.Item(33,3).Value = -17.73642719428278
$wsForecast.Cells.Item(33,4).Value = 98.01841555808389
$wsForecast.Cells.Item(34,1).Value = 45319.99999999999
$wsForecast.Cells.Item(34,2).Value = 40
$wsForecast.Cells.Item(34,3).Value = -16.56103399972706
$wsForecast.Cells.Item(34,4).Value = 99.04164219213277
$wsForecast.Cells.Item(35,1).Value = 45326.99999999999
$wsForecast.Cells.Item(35,2).Value = 40
$wsForecast.Cells.Item(35,3).Value = -16.10538449579345
$wsForecast.Cells.Item(35,4).Value = 100.9062262639273
$wsForecast.Cells.Item(36,1).Value = 45333.99999999999
$wsForecast.Cells.Item(36,2).Value = 41
$wsForecast.Cells.Item(36,3).Value = -16.45721202047351
$wsForecast.Cells.Item(36,4).Value = 100.055679302632

# 5. Apply the same date/time display format used on the other sheets' date column
$wsForecast.Range("A2:A36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
